$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 12: "CNN" -> "CNN " (trailing space) and fill in its metrics
$ws.Range("B12").Value = "CNN "
$ws.Range("C12").Value = 0.632
$ws.Range("D12").Value = 0.646
$ws.Range("E12").Value = 0.713
$ws.Range("F12").Value = 0.678

# Row 11: "ANN (1st layer = 10, 2nd Layer = 6)" -> "ANN " (trailing space preserved)
$ws.Range("B11").Value = "ANN "

# New row 16: CatBoost model results
$ws.Range("A16").Value = 2
$ws.Range("B16").Value = "CatBoost"
$ws.Range("C16").Value = 0.696
$ws.Range("D16").Value = 0.853
$ws.Range("E16").Value = 0.703
$ws.Range("F16").Value = 0.771

# Update selection to match target state
$ws.Range("A17").Select()
